# Update the "想去人数" (F) and "最低票价" (G) figures on the sheets that
# contain the exhibition listing data: "展览" and "全部类型" (both sheets
# hold a mirrored copy of the same table).
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 545
    $ws.Range("G2").Value = 75

    $ws.Range("F7").Value = 780

    $ws.Range("F9").Value = 417
}
